$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8410.643
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 8826.846
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 8826.846
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -9176.846
$ws.Range("H76").Value = 9332.666999999999
$ws.Range("I76").Value = 9995
$ws.Range("J76").Value = 9001.5
$ws.Range("K76").Value = 9995
$ws.Range("L76").Value = 9001.5
$ws.Range("M76").Value = -9680
$ws.Range("N76").Value = -9631.5
$ws.Range("H79").Value = 9332.666999999999
$ws.Range("I79").Value = 9995
$ws.Range("J79").Value = 9001.5
$ws.Range("K79").Value = 9995
$ws.Range("L79").Value = 9001.5
$ws.Range("M79").Value = -8903
$ws.Range("N79").Value = -11185.5
$ws.Range("H111").Value = 1876.2778
$ws.Range("I111").Value = 2158.8
$ws.Range("J111").Value = 1523.125
$ws.Range("K111").Value = 6476.400000000001
$ws.Range("L111").Value = 4569.375
$ws.Range("M111").Value = -3409.400000000001
$ws.Range("N111").Value = -10703.375
$ws.Range("H112").Value = 2739.6365
$ws.Range("J112").Value = 3187
$ws.Range("L112").Value = 9561
$ws.Range("N112").Value = -11777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14636
$ws.Range("I32").Value = 13693.818
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 13693.818
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -13406.818
$ws.Range("N32").Value = -25574
$ws.Range("H35").Value = 2784
$ws.Range("I35").Value = 2784
$ws.Range("K35").Value = 2784
$ws.Range("M35").Value = -2378
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2889.3809
$ws.Range("I20").Value = 2063.2856
$ws.Range("J20").Value = 4541.5713
$ws.Range("K20").Value = 2063.2856
$ws.Range("L20").Value = 4541.5713
$ws.Range("M20").Value = -1816.2856
$ws.Range("N20").Value = -5035.5713
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H99").Value = 1614.1428
$ws.Range("I99").Value = 1549.8334
$ws.Range("K99").Value = 1549.8334
$ws.Range("M99").Value = -51.83339999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1135.5454
$ws.Range("I22").Value = 1110
$ws.Range("K22").Value = 1110
$ws.Range("M22").Value = -760
$ws.Range("H31").Value = 3659
$ws.Range("I31").Value = 3238.8333
$ws.Range("K31").Value = 3238.8333
$ws.Range("M31").Value = -2943.8333
$ws.Range("H34").Value = 3659
$ws.Range("I34").Value = 3238.8333
$ws.Range("K34").Value = 3238.8333
$ws.Range("M34").Value = -3036.8333
$ws.Range("H56").Value = 6000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 6000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 6000
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -7690
$ws.Range("H94").Value = 3882.8572
$ws.Range("I94").Value = 4337
$ws.Range("J94").Value = 3542.25
$ws.Range("K94").Value = 4337
$ws.Range("L94").Value = 3542.25
$ws.Range("M94").Value = -3886
$ws.Range("N94").Value = -4444.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2157.6
$ws.Range("J131").Value = 2329.6667
$ws.Range("L131").Value = 6989.000100000001
$ws.Range("N131").Value = -17069.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2391.6667
$ws.Range("I80").Value = 2300
$ws.Range("J80").Value = 2483.3333
$ws.Range("K80").Value = 2300
$ws.Range("L80").Value = 2483.3333
$ws.Range("M80").Value = -1302
$ws.Range("N80").Value = -4479.3333
$ws.Range("H83").Value = 2391.6667
$ws.Range("I83").Value = 2300
$ws.Range("J83").Value = 2483.3333
$ws.Range("K83").Value = 11500
$ws.Range("L83").Value = 12416.6665
$ws.Range("M83").Value = -6508
$ws.Range("N83").Value = -22400.6665
$ws.Range("H97").Value = 1969
$ws.Range("I97").Value = 954.5
$ws.Range("K97").Value = 954.5
$ws.Range("M97").Value = -458.5
$ws.Range("H107").Value = 850
$ws.Range("I107").Value = 791.6667
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 791.6667
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 1128.3333
$ws.Range("N107").Value = -5040
$ws.Range("H126").Value = 124831600
$ws.Range("J126").Value = 4505.2
$ws.Range("L126").Value = 13515.6
$ws.Range("N126").Value = -18455.6
$ws.Range("H132").Value = 3980.5625
$ws.Range("I132").Value = 3828.4285
$ws.Range("J132").Value = 4098.8887
$ws.Range("K132").Value = 11485.2855
$ws.Range("L132").Value = 12296.6661
$ws.Range("M132").Value = -8955.2855
$ws.Range("N132").Value = -17356.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 45000
$ws.Range("J3").Value = 45000
$ws.Range("L3").Value = 45000
$ws.Range("N3").Value = -45224
$ws.Range("H4").Value = 49800
$ws.Range("J4").Value = 49800
$ws.Range("L4").Value = 49800
$ws.Range("N4").Value = -50026
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H10").Value = 1599.7142
$ws.Range("J10").Value = 2222.25
$ws.Range("L10").Value = 2222.25
$ws.Range("N10").Value = -2502.25
$ws.Range("H12").Value = 2940.8
$ws.Range("J12").Value = 3550.25
$ws.Range("L12").Value = 3550.25
$ws.Range("N12").Value = -3890.25
$ws.Range("H13").Value = 1431.8334
$ws.Range("I13").Value = 1026.25
$ws.Range("J13").Value = 2243
$ws.Range("K13").Value = 1026.25
$ws.Range("L13").Value = 2243
$ws.Range("M13").Value = -886.25
$ws.Range("N13").Value = -2523
$ws.Range("H15").Value = 45000
$ws.Range("J15").Value = 45000
$ws.Range("L15").Value = 45000
$ws.Range("N15").Value = -45340
$ws.Range("H17").Value = 14908
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 19710.666
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 19710.666
$ws.Range("M17").Value = -330
$ws.Range("N17").Value = -20050.666
$ws.Range("H19").Value = 6600.8
$ws.Range("J19").Value = 8200.25
$ws.Range("L19").Value = 8200.25
$ws.Range("N19").Value = -8540.25
$ws.Range("H20").Value = 3800
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 1400
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 1400
$ws.Range("M20").Value = -4774
$ws.Range("N20").Value = -1852
$ws.Range("H25").Value = 30983.334
$ws.Range("I25").Value = 32500
$ws.Range("J25").Value = 30225
$ws.Range("K25").Value = 32500
$ws.Range("L25").Value = 30225
$ws.Range("M25").Value = -32270
$ws.Range("N25").Value = -30685
$ws.Range("H28").Value = 49800
$ws.Range("J28").Value = 49800
$ws.Range("L28").Value = 49800
$ws.Range("N28").Value = -50264
$ws.Range("H30").Value = 480.375
$ws.Range("I30").Value = 480.375
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 480.375
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -372.375
$ws.Range("N30").ClearContents()
$ws.Range("H37").Value = 49800
$ws.Range("J37").Value = 49800
$ws.Range("L37").Value = 49800
$ws.Range("N37").Value = -50014
$ws.Range("H55").Value = 661.3333
$ws.Range("I55").Value = 315.55554
$ws.Range("K55").Value = 315.55554
$ws.Range("M55").Value = -142.55554
$ws.Range("H93").Value = 2667.3333
$ws.Range("I93").Value = 2667.3333
$ws.Range("K93").Value = 2667.3333
$ws.Range("M93").Value = -1419.3333
$ws.Range("H136").Value = 5500
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 1025.6
$ws.Range("I21").Value = 504.33334
$ws.Range("J21").Value = 5717
$ws.Range("K21").Value = 504.33334
$ws.Range("L21").Value = 5717
$ws.Range("M21").Value = -269.33334
$ws.Range("N21").Value = -6187
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H35").Value = 1025.6
$ws.Range("I35").Value = 504.33334
$ws.Range("J35").Value = 5717
$ws.Range("K35").Value = 504.33334
$ws.Range("L35").Value = 5717
$ws.Range("M35").Value = -214.33334
$ws.Range("N35").Value = -6297
$ws.Range("H75").Value = 24974.75
$ws.Range("I75").Value = 24949.5
$ws.Range("K75").Value = 24949.5
$ws.Range("M75").Value = -24013.5
$ws.Range("H78").Value = 24974.75
$ws.Range("I78").Value = 24949.5
$ws.Range("K78").Value = 74848.5
$ws.Range("M78").Value = -70168.5
$ws.Range("H81").Value = 1269
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 1269
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 15000
$ws.Range("N84").Value = -25608
$ws.Range("H107").Value = 1304.1666
$ws.Range("I107").Value = 1165.4
$ws.Range("K107").Value = 3496.2
$ws.Range("M107").Value = -1576.2
$ws.Range("H122").Value = 724.75
$ws.Range("J122").Value = 724.5
$ws.Range("L122").Value = 2173.5
$ws.Range("N122").Value = -7073.5
$ws.Range("H136").Value = 31719.47
$ws.Range("I136").Value = 37230.785
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 111692.355
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -109142.355
$ws.Range("N136").Value = -23100
